$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update date/time values in column A (rows 2-5) ---
$ws.Range("A2").Value = 45170.50694444445
$ws.Range("A3").Value = 45170.51388888889
$ws.Range("A4").Value = 45170.52083333334
$ws.Range("A5").Value = 45170.52777777778

# --- Update data matrix B2:AH5 ---
$data = New-Object 'object[,]' 4,33
$data[0,0] = 13.798
$data[0,1] = 9.137
$data[0,2] = 3.527
$data[0,3] = 29.879
$data[0,4] = 22.444
$data[0,5] = 10.657
$data[0,6] = 31.967
$data[0,7] = 16.87
$data[0,8] = 6.742
$data[0,9] = 10.011
$data[0,10] = 11.733
$data[0,11] = 12.516
$data[0,12] = 3.497
$data[0,13] = 10.903
$data[0,14] = 14.966
$data[0,15] = 9.704000000000001
$data[0,16] = 3.059
$data[0,17] = 1.672
$data[0,18] = 158.575
$data[0,19] = 30.182
$data[0,20] = 10.064
$data[0,21] = 19.331
$data[0,22] = 9.888999999999999
$data[0,23] = 2.87
$data[0,24] = 17.102
$data[0,25] = 8.888999999999999
$data[0,26] = 8.15
$data[0,27] = 9.673999999999999
$data[0,28] = 12.058
$data[0,29] = 3.072
$data[0,30] = 28.934
$data[0,31] = 5.407
$data[0,32] = 12.581
$data[1,0] = 14.782
$data[1,1] = 10.603
$data[1,2] = 1.706
$data[1,3] = 32.373
$data[1,4] = 25.55
$data[1,5] = 11.51
$data[1,6] = 43.925
$data[1,7] = 18.033
$data[1,8] = 7.745
$data[1,9] = 11.315
$data[1,10] = 12.902
$data[1,11] = 13.758
$data[1,12] = 3.744
$data[1,13] = 11.655
$data[1,14] = 16.372
$data[1,15] = 10.156
$data[1,16] = 1.437
$data[1,17] = 0.974
$data[1,18] = 170.073
$data[1,19] = 32.594
$data[1,20] = 10.758
$data[1,21] = 21.48
$data[1,22] = 11.198
$data[1,23] = 2.149
$data[1,24] = 21.811
$data[1,25] = 9.502000000000001
$data[1,26] = 8.589
$data[1,27] = 10.12
$data[1,28] = 13.385
$data[1,29] = 1.198
$data[1,30] = 40.343
$data[1,31] = 5.914
$data[1,32] = 13.45
$data[2,0] = 19.122
$data[2,1] = 14.076
$data[2,2] = 1.395
$data[2,3] = 41.815
$data[2,4] = 33.733
$data[2,5] = 14.964
$data[2,6] = 57.809
$data[2,7] = 23.269
$data[2,8] = 10.227
$data[2,9] = 15.052
$data[2,10] = 16.734
$data[2,11] = 17.804
$data[2,12] = 4.831
$data[2,13] = 15.038
$data[2,14] = 21.309
$data[2,15] = 12.827
$data[2,16] = 1.019
$data[2,17] = 0.894
$data[2,18] = 221.583
$data[2,19] = 42.046
$data[2,20] = 13.881
$data[2,21] = 28.076
$data[2,22] = 14.705
$data[2,23] = 2.382
$data[2,24] = 28.299
$data[2,25] = 12.261
$data[2,26] = 10.941
$data[2,27] = 12.872
$data[2,28] = 17.483
$data[2,29] = 0.746
$data[2,30] = 52.703
$data[2,31] = 7.751
$data[2,32] = 17.354
$data[3,0] = 20.09
$data[3,1] = 14.9
$data[3,2] = 1.22
$data[3,3] = 43.91
$data[3,4] = 35.66
$data[3,5] = 15.75
$data[3,6] = 61.85
$data[3,7] = 24.43
$data[3,8] = 10.81
$data[3,9] = 15.96
$data[3,10] = 17.59
$data[3,11] = 18.7
$data[3,12] = 5.07
$data[3,13] = 15.79
$data[3,14] = 22.44
$data[3,15] = 13.37
$data[3,16] = 0.8100000000000001
$data[3,17] = 0.83
$data[3,18] = 233.03
$data[3,19] = 44.17
$data[3,20] = 14.58
$data[3,21] = 29.62
$data[3,22] = 15.52
$data[3,23] = 2.38
$data[3,24] = 30.07
$data[3,25] = 12.87
$data[3,26] = 11.44
$data[3,27] = 13.45
$data[3,28] = 18.42
$data[3,29] = 0.54
$data[3,30] = 56.31
$data[3,31] = 8.18
$data[3,32] = 18.22
$ws.Range("B2:AH5").Value = $data

# --- Column width changes (7 -> 8 in xml terms); ColumnWidth = 7.17 gives xml width 8 ---
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(23).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17

# --- Remove row 6 (data trimmed to 1000 -> here reduced sample rows) ---
$ws.Rows.Item(6).Delete()

Write-Output "done"
